# Weekly update: insert two new price records (Primera / Segunda) for the
# most recent reporting date at the top of the Pomelo "Start Ruby" data
# block (rows 136-137), pushing all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 136 (shifts old rows
# 136..228 down to 138..230, matching the dimension growing to A1:T230).
$ws.Range("A136:A137").EntireRow.Insert()

# --- New row 136: Primera ---
$ws.Cells.Item(136, 1).Value = 4
$ws.Cells.Item(136, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(136, 3).Value = "Los Lagos"
$ws.Cells.Item(136, 4).Value = 44603
$ws.Cells.Item(136, 5).Value = 10
$ws.Cells.Item(136, 6).Value = "Fruta"
$ws.Cells.Item(136, 7).Value = 100102
$ws.Cells.Item(136, 8).Value = "Cítricos"
$ws.Cells.Item(136, 9).Value = 100102006
$ws.Cells.Item(136, 10).Value = "Pomelo"
$ws.Cells.Item(136, 11).Value = "Start Ruby"
$ws.Cells.Item(136, 12).Value = "Primera"
$ws.Cells.Item(136, 13).Value = 200
$ws.Cells.Item(136, 14).Value = 13000
$ws.Cells.Item(136, 15).Value = 14000
$ws.Cells.Item(136, 16).Value = 13500
$ws.Cells.Item(136, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(136, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(136, 19).Value = 964
$ws.Cells.Item(136, 20).Value = 14

# --- New row 137: Segunda ---
$ws.Cells.Item(137, 1).Value = 4
$ws.Cells.Item(137, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(137, 3).Value = "Los Lagos"
$ws.Cells.Item(137, 4).Value = 44603
$ws.Cells.Item(137, 5).Value = 10
$ws.Cells.Item(137, 6).Value = "Fruta"
$ws.Cells.Item(137, 7).Value = 100102
$ws.Cells.Item(137, 8).Value = "Cítricos"
$ws.Cells.Item(137, 9).Value = 100102006
$ws.Cells.Item(137, 10).Value = "Pomelo"
$ws.Cells.Item(137, 11).Value = "Start Ruby"
$ws.Cells.Item(137, 12).Value = "Segunda"
$ws.Cells.Item(137, 13).Value = 100
$ws.Cells.Item(137, 14).Value = 11000
$ws.Cells.Item(137, 15).Value = 11000
$ws.Cells.Item(137, 16).Value = 11000
$ws.Cells.Item(137, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(137, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(137, 19).Value = 786
$ws.Cells.Item(137, 20).Value = 14
